$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(9, 1).Value = (Get-Date -Year 2016 -Month 8 -Day 29 -Hour 21 -Minute 12 -Second 31)
$ws.Cells.Item(9, 1).NumberFormat = "m/d/yy h:mm"

$ws.Cells.Item(9, 2).Value = 8
$ws.Cells.Item(9, 3).Value = 51
$ws.Cells.Item(9, 4).Value = 47
$ws.Cells.Item(9, 5).Value = 83
$ws.Cells.Item(9, 6).Value = 16
$ws.Cells.Item(9, 7).Value = 19635
$ws.Cells.Item(9, 8).Value = 18974
$ws.Cells.Item(9, 9).Value = 1081
$ws.Cells.Item(9, 10).Value = 206
$ws.Cells.Item(9, 11).Value = 189
$ws.Cells.Item(9, 12).Value = 5
$ws.Cells.Item(9, 13).Value = 1
$ws.Cells.Item(9, 14).Value = "Named"
